$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("A19").Value = "Visualisation"
$ws.Range("B19").Value = 15

$ws.Range("B20").Select()
